$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.288.68"
$ws.Range("E2").Value = "  +0.63%  "

$ws.Range("D3").Value = "3.767.84"
$ws.Range("E3").Value = "  +0.48%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'621.32"
$ws.Range("E5").Value = "  +0.58%  "

$ws.Range("D6").Value = "'181.50"
$ws.Range("E6").Value = "  +2.10%  "

$ws.Range("D7").Value = "3.765.73"
$ws.Range("E7").Value = "  +0.44%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  -1.71%  "

$ws.Range("E10").Value = "  -1.05%  "

$ws.Range("D11").Value = "'6.56"
$ws.Range("E11").Value = "  +2.42%  "

$ws.Range("D12").Value = "'0.485"
$ws.Range("E12").Value = "  -3.56%  "

$ws.Range("D13").Value = "'40.19"
$ws.Range("E13").Value = "  -2.17%  "

$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").Value = "4.397.32"
$ws.Range("E15").Value = "  +0.46%  "

$ws.Range("D16").Value = "3.774.08"
$ws.Range("E16").Value = "  +0.65%  "

$ws.Range("D17").Value = "70.284.02"
$ws.Range("E17").Value = "  +0.42%  "

$ws.Range("D18").Value = "'7.59"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("E19").Value = "  -2.15%  "

$ws.Range("D20").Value = "'16.72"
$ws.Range("E20").Value = "  -0.40%  "

$ws.Range("D21").Value = "'507.22"
$ws.Range("E21").Value = "  -2.14%  "

$ws.Range("D22").Value = "'9.26"
$ws.Range("E22").Value = "  -1.40%  "

$ws.Range("E23").Value = "  -1.43%  "

$ws.Range("D24").Value = "'2.63"
$ws.Range("E24").Value = "  +5.28%  "

$ws.Range("D25").Value = "'87.03"
$ws.Range("E25").Value = "  -2.12%  "

$ws.Range("D26").Value = "'11.37"
$ws.Range("E26").Value = "  +3.81%  "

$ws.Range("E27").Value = "  -3.68%  "

$ws.Range("E28").Value = "  +6.80%  "

$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("E30").Value = "  -0.27%  "

$ws.Range("E31").Value = "  +2.66%  "

$ws.Range("E32").Value = "  +1.72%  "

$ws.Range("D33").Value = "'30.89"
$ws.Range("E33").Value = "  -2.36%  "

$ws.Range("D34").Value = "'0.114"
$ws.Range("E34").Value = "  -1.13%  "

$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.05%  "

$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").Value = "'1.06"
$ws.Range("E36").Value = "  +1.27%  "

$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").Value = "'6.21"
$ws.Range("E37").Value = "  -0.18%  "

$ws.Range("E38").Value = "  +2.95%  "

$ws.Range("E39").Value = "  +5.84%  "

$ws.Range("D40").Value = "'3.10"
$ws.Range("E40").Value = "  +13.31%  "

$ws.Range("D41").Value = "'2.08"
$ws.Range("E41").Value = "  -5.12%  "

$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").Value = "'442.72"
$ws.Range("E42").Value = "  +3.56%  "

$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D43").Value = "'45.94"
$ws.Range("E43").Value = "  +2.80%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "'49.93"
$ws.Range("E44").Value = "  -2.99%  "

$ws.Range("E45").Value = "  -2.13%  "

$ws.Range("D46").Value = "2.993.02"
$ws.Range("E46").Value = "  -2.81%  "

$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("D48").Value = "'27.68"
$ws.Range("E48").Value = "  -0.76%  "

$ws.Range("D49").Value = "'139.11"
$ws.Range("E49").Value = "  +1.95%  "

$ws.Range("E50").Value = "  -0.04%  "

$ws.Range("D51").Value = "'2.49"
$ws.Range("E51").Value = "  -1.05%  "
